# AVC_CmdAndTlmAndTest_v2.xlsx - "Rpi to Gui OccGrid" sheet rework
# Expands the Occupancy Grid telemetry packet layout: splits the old
# "Angle" field into the new Rows/Cols/Car-position/Angle-result fields,
# renames the sync word, widens columns A/B, drops the header's thick
# bottom border, and adds a trailing Checksum field as the new final row.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Rpi to Gui OccGrid")

# ---------------------------------------------------------------------
# 1. Make room: insert 4 blank rows above the old row 5 (the
#    "Occcupancy Grid" row), pushing it down to row 9, the blank
#    row after it to row 10, and the trailing blank row to row 11.
# ---------------------------------------------------------------------
$ws.Rows.Item(5).Insert()
$ws.Rows.Item(5).Insert()
$ws.Rows.Item(5).Insert()
$ws.Rows.Item(5).Insert()

# ---------------------------------------------------------------------
# 2. Row 2: sync word renamed
# ---------------------------------------------------------------------
$ws.Range("B2").Value = "Always 0x55555555"

# ---------------------------------------------------------------------
# 3. Row 4 (was "Angle") becomes "Number of Rows"
# ---------------------------------------------------------------------
$ws.Range("A4").Value = "Number of Rows"
$ws.Range("B4").Value = "signed integer"
$ws.Range("C4").Value = 1
$ws.Range("D4").Value = 4
$ws.Range("E4").Formula = "=C4*D4"
$ws.Range("F4").Formula = "=F3+E4"

# Rows 5-7: Number of Cols / Car position Row / Car Position Col
$ws.Range("A5").Value = "Number of Cols"
$ws.Range("B5").Value = "signed integer"
$ws.Range("C5").Value = 1
$ws.Range("D5").Value = 4
$ws.Range("E5").Formula = "=C5*D5"
$ws.Range("F5").Formula = "=F4+E5"

$ws.Range("A6").Value = "Car position Row"
$ws.Range("B6").Value = "signed integer"
$ws.Range("C6").Value = 1
$ws.Range("D6").Value = 4
$ws.Range("E6").Formula = "=C6*D6"
$ws.Range("F6").Formula = "=F5+E6"

$ws.Range("A7").Value = "Car Position Col"
$ws.Range("B7").Value = "signed integer"
$ws.Range("C7").Value = 1
$ws.Range("D7").Value = 4
$ws.Range("E7").Formula = "=C7*D7"
$ws.Range("F7").Formula = "=F6+E7"

# Row 8: Angle result (two-line description -> wraps)
$ws.Range("A8").Value = "Angle result"
$ws.Range("B8").Value = "Angle of analysis result - signed int.  Positive is CW neg is CCW"
$ws.Range("C8").Value = 1
$ws.Range("D8").Value = 4
$ws.Range("E8").Formula = "=C8*D8"
$ws.Range("F8").Formula = "=F7+E8"

# ---------------------------------------------------------------------
# 4. Row 9 (the old "Occcupancy Grid" row, shifted down by the insert)
#    is no longer the last row, so its content gets the normal
#    (non-thick-bottom) row style, and the description + formula are
#    refreshed.
# ---------------------------------------------------------------------
$ws.Range("A9").Value = "Occcupancy Grid"
$ws.Range("B9").Value = "Array of bytes (300 x 50 TBD)"
$ws.Range("C9").Value = 15000
$ws.Range("D9").Value = 1
$ws.Range("E9").Formula = "=C9*D9"
$ws.Range("F9").Formula = "=F8+E9"

$ws.Range("A9:F9").Borders.Item(9).LineStyle = 0
$ws.Range("A9:F9").Borders.Item(7).LineStyle = 0
$ws.Range("A9:F9").Borders.Item(10).LineStyle = 0
$ws.Range("A9").Borders.Item(7).LineStyle = 1
$ws.Range("A9").Borders.Item(7).Weight = -4138
$ws.Range("F9").Borders.Item(10).LineStyle = 1
$ws.Range("F9").Borders.Item(10).Weight = -4138
$ws.Range("A9:F9").Font.Bold = $false

# ---------------------------------------------------------------------
# 5. Row 10 (the old trailing blank row, shifted down) becomes the new
#    final row, with the "Checksum" field and the thick bottom border.
# ---------------------------------------------------------------------
$ws.Range("A10").Value = "Checksum"
$ws.Range("B10").Value = "Summation of all integer except for the checksum itself.  "
$ws.Range("C10").Value = 1
$ws.Range("D10").Value = 4
$ws.Range("E10").Formula = "=C10*D10"
$ws.Range("F10").Formula = "=F9+E10"

$ws.Range("A10:F10").Borders.Item(9).LineStyle = 1
$ws.Range("A10:F10").Borders.Item(9).Weight = -4138
$ws.Range("A10").Borders.Item(7).LineStyle = 1
$ws.Range("A10").Borders.Item(7).Weight = -4138
$ws.Range("F10").Borders.Item(10).LineStyle = 1
$ws.Range("F10").Borders.Item(10).Weight = -4138

# ---------------------------------------------------------------------
# 6. Header row: drop the thick bottom border (row no longer ht=73/
#    thickBot), widen columns A & B.
# ---------------------------------------------------------------------
$ws.Range("A1:F1").Borders.Item(9).LineStyle = 0
$ws.Columns.Item(1).ColumnWidth = 23.36328125
$ws.Columns.Item(2).ColumnWidth = 30.08984375

# ---------------------------------------------------------------------
# 7. Row heights: recompute for the new layout.
# ---------------------------------------------------------------------
$ws.Rows.Item(1).RowHeight = 58
$ws.Rows.Item(2).RowHeight = 14.5
$ws.Rows.Item(3).RowHeight = 14.5
$ws.Rows.Item(4).RowHeight = 14.5
$ws.Rows.Item(5).RowHeight = 14.5
$ws.Rows.Item(6).RowHeight = 14.5
$ws.Rows.Item(7).RowHeight = 14.5
$ws.Rows.Item(8).RowHeight = 29
$ws.Rows.Item(9).RowHeight = 14.5
$ws.Rows.Item(10).RowHeight = 29.5

# ---------------------------------------------------------------------
# 8. Selection cosmetics, dimension naturally grows to A1:I11.
# ---------------------------------------------------------------------
$ws.Range("H12").Select()
